# "Assign Task Test Case"
# Simulate the sequence of task assignments recorded while running the
# automated Katalon test, ending with the final assigned task timestamp
# that is committed into cell B2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "AUTO_TOOL_KAT_TASK_02/11/2018-10:47:56"
$ws.Range("B2").Value = "AUTO_TOOL_KAT_TASK_02/11/2018-11:19:37"
$ws.Range("B2").Value = "AUTO_TOOL_KAT_TASK_02/11/2018-11:24:17"
$ws.Range("B2").Value = "AUTO_TOOL_KAT_TASK_02/11/2018-11:29:11"
$ws.Range("B2").Value = "AUTO_TOOL_KAT_TASK_02/11/2018-11:31:25"
$ws.Range("B2").Value = "AUTO_TOOL_KAT_TASK_02/11/2018-11:37:10"
